$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.000699520111084
$ws.Range("B1").Value = 2.111485719680786
$ws.Range("C1").Value = 6.855403423309326
$ws.Range("D1").Value = 1.998162150382996
$ws.Range("E1").Value = 1.373103618621826
